$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.107.38'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.878.50'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.07%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.54'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.11%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5042'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3958'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08230'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.72%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.15'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.86%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.094'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.80%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.876.16'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.299'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.200'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.72%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.76'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.77%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001088'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06466'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.36%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.09%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '30.087.63'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.842'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.93%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.75%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.170'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.097.62'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.37%  '

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.20'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.65%  '

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.85'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.59%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.249'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -7.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.44'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.073'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1035'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.939'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.700'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02424'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.296'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06359'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.83%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2134'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.174'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.505'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6302'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.216'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.73%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.28%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.19'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.82%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5907'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.094'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.71%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.627'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.209'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.20'

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.122'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.88%  '
